$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.8966636666666666
$ws.Range("H2").Value = 2.689991
$ws.Range("I2").Value = 0.334725143386341
$ws.Range("J2").Value = 0.3647360854412732
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2416746666666667
$ws.Range("N2").Value = 0.725024
$ws.Range("O2").Value = 0.08321776967690767
$ws.Range("P2").Value = 0.08321776967690767
$ws.Range("Q2").Value = 0.2167008927537778
$ws.Range("R2").Value = 1.950308034784
$ws.Range("S2").Value = 0.02785507988739442
$ws.Range("T2").Value = 0.03035252355110879
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.8966636666666666
$ws.Range("H3").Value = 2.689991
$ws.Range("I3").Value = 0.334725143386341
$ws.Range("J3").Value = 0.3647360854412732
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.662448666666667
$ws.Range("N3").Value = 7.987346000000001
$ws.Range("O3").Value = 0.9167822303230924
$ws.Range("P3").Value = 0.9167822303230924
$ws.Range("Q3").Value = 2.387320983765111
$ws.Range("R3").Value = 21.485888853886
$ws.Range("S3").Value = 0.3068700634989466
$ws.Range("T3").Value = 0.3343835618901644
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 0.6827986666666667
$ws.Range("I4").Value = 0.2548891965854188
$ws.Range("J4").Value = 0.2777421703171357
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.2416746666666667
$ws.Range("N4").Value = 0.725024
$ws.Range("O4").Value = 0.08321776967690767
$ws.Range("P4").Value = 0.08321776967690767
$ws.Range("Q4").Value = 0.1650151401671111
$ws.Range("R4").Value = 1.485136261504
$ws.Range("S4").Value = 0.02121131045457742
$ws.Range("T4").Value = 0.02311308395901586
$ws.Range("D5").Value = "FAPs"
$ws.Range("G5").Value = 0.6827986666666667
$ws.Range("I5").Value = 0.2548891965854188
$ws.Range("J5").Value = 0.2777421703171357
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.662448666666667
$ws.Range("N5").Value = 7.987346000000001
$ws.Range("O5").Value = 0.9167822303230924
$ws.Range("P5").Value = 0.9167822303230924
$ws.Range("Q5").Value = 1.817916399668444
$ws.Range("R5").Value = 16.361247597016
$ws.Range("S5").Value = 0.2336778861308414
$ws.Range("T5").Value = 0.2546290863581199
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.221369
$ws.Range("H6").Value = 0.664107
$ws.Range("I6").Value = 0.08263719499391366
$ws.Range("J6").Value = 0.09004631892602898
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.2416746666666667
$ws.Range("N6").Value = 0.725024
$ws.Range("O6").Value = 0.08321776967690767
$ws.Range("P6").Value = 0.08321776967690767
$ws.Range("Q6").Value = 0.05349927928533334
$ws.Range("R6").Value = 0.481493513568
$ws.Range("S6").Value = 0.006876883059749214
$ws.Range("T6").Value = 0.007493453828639651
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.221369
$ws.Range("H7").Value = 0.664107
$ws.Range("I7").Value = 0.08263719499391366
$ws.Range("J7").Value = 0.09004631892602898
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.662448666666667
$ws.Range("N7").Value = 7.987346000000001
$ws.Range("O7").Value = 0.9167822303230924
$ws.Range("P7").Value = 0.9167822303230924
$ws.Range("Q7").Value = 0.5893835988913334
$ws.Range("R7").Value = 5.304452390022
$ws.Range("S7").Value = 0.07576031193416445
$ws.Range("T7").Value = 0.08255286509738932
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 0.6612465000000001
$ws.Range("H8").Value = 1.322493
$ws.Range("I8").Value = 0.2468437584284291
$ws.Range("J8").Value = 0.1793169270244717
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.2416746666666667
$ws.Range("N8").Value = 0.725024
$ws.Range("O8").Value = 0.08321776967690767
$ws.Range("P8").Value = 0.08321776967690767
$ws.Range("Q8").Value = 0.159806527472
$ws.Range("R8").Value = 0.9588391648320002
$ws.Range("S8").Value = 0.02054178703507925
$ws.Range("T8").Value = 0.01492235473229335
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 0.6612465000000001
$ws.Range("H9").Value = 1.322493
$ws.Range("I9").Value = 0.2468437584284291
$ws.Range("J9").Value = 0.1793169270244717
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.662448666666667
$ws.Range("N9").Value = 7.987346000000001
$ws.Range("O9").Value = 0.9167822303230924
$ws.Range("P9").Value = 0.9167822303230924
$ws.Range("Q9").Value = 1.760534862263
$ws.Range("R9").Value = 10.563209173578
$ws.Range("S9").Value = 0.2263019713933499
$ws.Range("T9").Value = 0.1643945722921784
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.216728
$ws.Range("H10").Value = 0.650184
$ws.Range("I10").Value = 0.08090470660589748
$ws.Range("J10").Value = 0.08815849829109049
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.2416746666666667
$ws.Range("N10").Value = 0.725024
$ws.Range("O10").Value = 0.08321776967690767
$ws.Range("P10").Value = 0.08321776967690767
$ws.Range("Q10").Value = 0.05237766715733334
$ws.Range("R10").Value = 0.471399004416
$ws.Range("S10").Value = 0.006732709240107366
$ws.Range("T10").Value = 0.007336353605850026
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.216728
$ws.Range("H11").Value = 0.650184
$ws.Range("I11").Value = 0.08090470660589748
$ws.Range("J11").Value = 0.08815849829109049
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.662448666666667
$ws.Range("N11").Value = 7.987346000000001
$ws.Range("O11").Value = 0.9167822303230924
$ws.Range("P11").Value = 0.9167822303230924
$ws.Range("Q11").Value = 0.5770271746293334
$ws.Range("R11").Value = 5.193244571664001
$ws.Range("S11").Value = 0.07417199736579011
$ws.Range("T11").Value = 0.08082214468524046
